# Add a new "2022-Q1" sheet (between "2021-Q4" and "总计"), populate it with
# fund-holding data, and prepend a corresponding summary row to "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q1" worksheet right after "2021-Q4" (i.e. right
#    before "总计", which is currently the last tab).
# ---------------------------------------------------------------------------
$q4Sheet   = $wb.Worksheets.Item(5)
$newSheet  = $wb.Worksheets.Add($null, $q4Sheet)
$newSheet.Name = "2022-Q1"

# NOTE: fetch the "总计" sheet by name (not by a stale index/object captured
# before the Add() above) since inserting a new tab shifts what a
# previously-held reference resolves to.
$totalSheet = $wb.Worksheets.Item("总计")

# Reference formatted cells on an existing quarter sheet ("2021-Q4") so we can
# clone their exact cell style (bold/border/center header style = style index
# 2, and the bold/border/center index-column style also = style index 2) onto
# the new sheet without introducing brand-new style entries.
$headerStyleSrc = $q4Sheet.Range("B1")
$indexStyleSrc  = $q4Sheet.Range("A2")

# ---------------------------------------------------------------------------
# Header row
# ---------------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = [char]([int][char]'B' + $i)
    $cell = $newSheet.Range("$col`1")
    $headerStyleSrc.Copy()
    $cell.PasteSpecial(-4122) # xlPasteFormats
    $cell.Value = $headers[$i]
}

# ---------------------------------------------------------------------------
# Data rows: index, code, name, scale, position, pct, value, rank
# ---------------------------------------------------------------------------
$rows = @(
    @(0,  "011300", "易方达智造优势混合A",                         "63.40", "89.88", "4.75", "3.0115", 10),
    @(1,  "011822", "易方达产业升级一年封闭运作混合型证券投资基金A", "65.20", "64.75", "4.60", "2.9992", 7),
    @(2,  "009049", "易方达高端制造混合",                           "52.81", "90.84", "4.63", "2.4451", 10),
    @(3,  "213001", "宝盈鸿利收益灵活配置混合A",                    "17.98", "90.37", "7.79", "1.4006", 3),
    @(4,  "012301", "易方达核心智造混合",                           "22.83", "61.70", "4.74", "1.0821", 6),
    @(5,  "011301", "易方达智造优势混合C",                          "9.96",  "89.88", "4.75", "0.4731", 10),
    @(6,  "011823", "易方达产业升级一年封闭运作混合型证券投资基金C", "10.22", "64.75", "4.60", "0.4701", 7),
    @(7,  "010751", "宝盈优质成长混合A",                            "5.64",  "92.80", "6.16", "0.3474", 4),
    @(8,  "001543", "宝盈新锐灵活配置混合A",                        "3.21",  "93.26", "7.51", "0.2411", 2),
    @(9,  "206002", "鹏华精选成长混合",                             "4.48",  "92.68", "4.74", "0.2124", 6),
    @(10, "233009", "大摩多因子精选策略混合",                       "6.77",  "89.73", "1.75", "0.1185", 2),
    @(11, "007581", "宝盈鸿利收益灵活配置混合C",                    "0.73",  "90.37", "7.79", "0.0569", 3),
    @(12, "010752", "宝盈优质成长混合C",                            "0.78",  "92.80", "6.16", "0.0480", 4),
    @(13, "007578", "宝盈新锐灵活配置混合C",                        "0.20",  "93.26", "7.51", "0.0150", 2),
    @(14, "009918", "上银核心成长混合A",                            "0.13",  "91.71", "0.79", "0.0010", 10),
    @(15, "009919", "上银核心成长混合C",                            "0.07",  "91.71", "0.79", "0.0006", 10)
)

$rowNum = 2
foreach ($r in $rows) {
    $indexStyleSrc.Copy()
    $newSheet.Range("A$rowNum").PasteSpecial(-4122)
    $newSheet.Range("A$rowNum").Value = $r[0]

    # Text-valued columns: force text storage (NumberFormat "@") so numeric
    # looking strings ("011300", "63.40", ...) keep their exact literal
    # representation (leading zeros / trailing zeros) instead of being
    # coerced into numbers, then restore the default "Normal" style so no
    # visible formatting change is left behind on the cell.
    $cols = @("B", "C", "D", "E", "F", "G")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cell = $newSheet.Range("$($cols[$i])$rowNum")
        $cell.NumberFormat = "@"
        $cell.Value = $r[$i + 1]
        $cell.Style = "Normal"
    }

    $newSheet.Range("H$rowNum").Value = $r[7]

    $rowNum++
}

# ---------------------------------------------------------------------------
# 2. Prepend a "2022-Q1" row to the "总计" summary sheet, shifting the
#    existing rows down by one and renumbering the index column.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert(-4121) # xlShiftDown

$indexStyleSrc2 = $totalSheet.Range("A3")
$indexStyleSrc2.Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 16
$totalSheet.Range("D2").Value = 12.92
$totalSheet.Range("B2:D2").Style = "Normal"

# The row-insert pushed the pre-existing rows down a slot but their index
# (col A) values are still the OLD 0..4 sequence; renumber them to 1..5 so
# the running index stays contiguous under the new row 0.
for ($r = 3; $r -le 7; $r++) {
    $totalSheet.Range("A$r").Value = $r - 2
}
